$d = $word.ActiveDocument

# --- Step 1: Insert a brand-new "NFR2" paragraph right before the "NFR3" paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "NFR3:*") {
        $target = $p.Range
        $target.Collapse(1)          # 1 = wdCollapseStart
        $target.InsertParagraphBefore()
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "NFR3:*") {
        $nfr3Para = $p
        break
    }
}
$newPara = $nfr3Para.Previous()
$newRange = $newPara.Range
$startPos = $newRange.Start
$newRange.InsertBefore("NFR2")

# format only the "NFR2" label text as bold (must not include the paragraph
# mark itself, otherwise the bold would leak into the paragraph's pPr/rPr)
$labelRange = $d.Range($startPos, $startPos + 4)
$labelRange.Bold = 1

# append the (non-bold) description text right after the label
$descPoint = $d.Range($startPos + 4, $startPos + 4)
$descPoint.InsertAfter(": The system shall be accessible through a graphical interface, with all UI functionality discoverable and accessible.")
$descPoint.Bold = 0

# --- Step 2: Replace NFR3's description text with the new NFR3 text
#     (this is the (slightly reworded) requirement that used to live in NFR5) ---
$d.Content.Find.Execute(
    ": The system shall be accessible through a graphical interface.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ":The system must let users set a primary and secondary color through an in-app setting, applying changes immediately or on restart without needing a reinstall. The default colors shall be UVU dark green and white.",
    2)

# --- Step 3: Delete the NFR4 paragraph entirely ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "NFR4:*") {
        $p.Range.Delete()
        break
    }
}

# --- Step 4: Delete the (old) NFR5 paragraph entirely ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "NFR5:*") {
        $p.Range.Delete()
        break
    }
}
